$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("441").Delete()
